$d = $word.ActiveDocument

# Replace the whole document story with the updated body content.
# This avoids leftover stray <w:proofErr/> markers from the original
# runs (Word keeps those attached to the paragraph independent of the
# text range contents), and lets us fully control run/proofErr
# structure for the rewritten first paragraph as well as insert the
# new empty paragraph before the final one.

$newBodyXml = '<w:p w:rsidR="00D21123" w:rsidRDefault="00D707EE"><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Mmm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> parece que </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>ste</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> se m demora </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>actua&#241;lizar</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' + `
    '<w:p w:rsidR="006F06BB" w:rsidRDefault="006F06BB"><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
    '<w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p>' + `
    '<w:p w:rsidR="00661E5F" w:rsidRPr="00D707EE" w:rsidRDefault="00661E5F"><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $newBodyXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xml)
